$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.965.24'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '1.634.27'
$ws.Range("E3").Value = '  -0.53%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.65%  '
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.45'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.71%  '
$ws.Range("E9").Value = '  -2.03%  '
$ws.Range("E10").Value = '  -0.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0881'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").Value = '1.628.64'
$ws.Range("E13").Value = '  -1.13%  '
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("E15").Value = '  -2.02%  '
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("D17").Value = '27.967.24'
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.66%  '
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("E22").Value = '  -6.24%  '
$ws.Range("E23").Value = '  -0.83%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.52'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  -0.37%  '
$ws.Range("E28").Value = '  -0.78%  '
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("E30").Value = '  -0.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0481'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.01%  '
$ws.Range("E32").Value = '  +1.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.10'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.35%  '
$ws.Range("D34").Value = '1.408.23'
$ws.Range("E34").Value = '  -1.15%  '
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("E36").Value = '  +10.78%  '
$ws.Range("E37").Value = '  +0.63%  '
$ws.Range("E38").Value = '  +1.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.556'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.867'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.61%  '
$ws.Range("E41").Value = '  -0.99%  '
$ws.Range("E42").Value = '  -0.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '66.83'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.58%  '
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("E45").Value = '  +0.65%  '
$ws.Range("E46").Value = '  -0.42%  '
$ws.Range("D47").Value = '1.775.37'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.16'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.92%  '
$ws.Range("D49").Value = '0.0₆0103'
$ws.Range("E49").Value = '  -2.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0999'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0505'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.34%  '
